# Edit: Tambah TC002-001 buat test kawalpemilu dan benerin document summary
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Update URL (E2) and PROJECT_NAME (K2) cell values to reflect new target site
$ws.Range("E2").Value = "https://kawalpemilu.org/"
$ws.Range("K2").Value = "Kawal Pemilu - Regression Test"

# Widen column K to fit new, longer text (stored width ends up as 29)
$ws.Columns.Item(11).ColumnWidth = 28.17

# Update the sheet's active cell / selection
$ws.Range("L7").Select()
